$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '  -3.62%  '
$ws.Range("D2").Value = '27.033.22'

$ws.Range("E3").Value = '  -3.88%  '
$ws.Range("D3").Value = '1.712.71'

$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"

$ws.Range("E5").Value = '  -5.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.56'
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = '  +0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  +5.84%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4758'
$ws.Range("D7").Style = "Normal"

$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3460'
$ws.Range("D8").Style = "Normal"

$ws.Range("E9").Value = '  +0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.34'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -2.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07297'
$ws.Range("D10").Style = "Normal"

$ws.Range("E11").Value = '  -5.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.042'
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"

$ws.Range("E13").Value = '  -5.37%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.82'
$ws.Range("D13").Style = "Normal"

$ws.Range("E14").Value = '  -3.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.848'
$ws.Range("D14").Style = "Normal"

$ws.Range("E15").Value = '  -3.74%  '
$ws.Range("D15").Value = '1.716.91'

$ws.Range("E16").Value = '  -6.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.811'
$ws.Range("D16").Style = "Normal"

$ws.Range("E17").Value = '  -4.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '89.08'
$ws.Range("D17").Style = "Normal"

$ws.Range("E18").Value = '  -2.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001038'
$ws.Range("D18").Style = "Normal"

$ws.Range("E19").Value = '  -0.82%  '

$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"

$ws.Range("E21").Value = '  -4.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.45'
$ws.Range("D21").Style = "Normal"

$ws.Range("E22").Value = '  -3.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.624'
$ws.Range("D22").Style = "Normal"

$ws.Range("E23").Value = '  -3.39%  '
$ws.Range("D23").Value = '27.138.51'

$ws.Range("E24").Value = '  -5.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.77'
$ws.Range("D24").Style = "Normal"

$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.094'
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = '  -6.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.78'
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = '  -2.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.81'
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = '  -3.93%  '
$ws.Range("D28").Value = '1.910.13'

$ws.Range("E29").Value = '  -4.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.063'
$ws.Range("D29").Style = "Normal"

$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '121.48'
$ws.Range("D30").Style = "Normal"

$ws.Range("E31").Value = '  -9.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.004'
$ws.Range("D31").Style = "Normal"

$ws.Range("E32").Value = '  +0.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09266'
$ws.Range("D32").Style = "Normal"

$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.605'
$ws.Range("D33").Style = "Normal"

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.341'
$ws.Range("D34").Style = "Normal"

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("E35").Value = '  +7.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.489'
$ws.Range("D35").Style = "Normal"

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E36").Value = '  -4.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02188'
$ws.Range("D36").Style = "Normal"

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("E37").Value = '  -4.86%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05855'
$ws.Range("D37").Style = "Normal"

$ws.Range("B38").Value = 'Aptos'
$ws.Range("C38").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("E38").Value = '  -7.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.04'
$ws.Range("D38").Style = "Normal"

$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E39").Value = '  -5.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1995'
$ws.Range("D39").Style = "Normal"

$ws.Range("B40").Value = 'Frax'
$ws.Range("C40").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.001'
$ws.Range("D40").Style = "Normal"

$ws.Range("B41").Value = 'InternetComputer(DFINITY)'
$ws.Range("C41").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("E41").Value = '  -5.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.728'
$ws.Range("D41").Style = "Normal"

$ws.Range("E42").Value = '  -6.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5934'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  -7.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.089'
$ws.Range("D43").Style = "Normal"

$ws.Range("E44").Value = '  -5.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.490'
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = '  -3.91%  '

$ws.Range("E46").Value = '  -4.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.587'
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = '  -5.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5570'
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = '  -3.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.37'
$ws.Range("D48").Style = "Normal"

$ws.Range("E49").Value = '  -6.46%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.832'
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  -2.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06701'
$ws.Range("D50").Style = "Normal"

$ws.Range("E51").Value = '  -5.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.083'
$ws.Range("D51").Style = "Normal"
